$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new trade row (row 9) - mirrors the existing repeater rows
$ws.Range("A9").Value = 42654.746261574073
$ws.Range("A9").NumberFormat = "m/d/yy h:mm"
$ws.Range("B9").Value = $true
$ws.Range("C9").Value = 9857.7199999999993
$ws.Range("D9").Value = 9850.33
$ws.Range("E9").Value = 308
$ws.Range("F9").Value = 308.45999999999998
$ws.Range("G9").Value = $false
$ws.Range("G9").NumberFormat = "m/d/yy h:mm"
$ws.Range("H9").Value = 0.15
$ws.Range("I9").Value = $false

# Re-run "best fit" column sizing on the used columns, like the repeater does
# after appending a new trade row (the data backing the auto-fit widths changed
# since the longest value per column is now possibly on row 9).
$ws.Columns("A").ColumnWidth = 14.5
$ws.Columns("B").ColumnWidth = 7.5
$ws.Columns("C").ColumnWidth = 8
$ws.Columns("D").ColumnWidth = 10.5
$ws.Columns("E").ColumnWidth = 10
$ws.Columns("F").ColumnWidth = 6.1666666666666667
$ws.Columns("G").ColumnWidth = 9.5
$ws.Columns("H").ColumnWidth = 13.8333333333333333
$ws.Columns("I").ColumnWidth = 11
